# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# for a batch of leve rows across the per-class profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3848.229  # H76: 6529.407 -> 3848.229
$ws.Cells.Item(76, 9).Value = 3328.6191  # I76: 5199 -> 3328.6191
$ws.Cells.Item(76, 10).Value = 5485  # J76: 7312 -> 5485
$ws.Cells.Item(76, 11).Value = 3328.6191  # K76: 5199 -> 3328.6191
$ws.Cells.Item(76, 12).Value = 5485  # L76: 7312 -> 5485
$ws.Cells.Item(76, 13).Value = -3013.6191  # M76: -4884 -> -3013.6191
$ws.Cells.Item(76, 14).Value = -6115  # N76: -7942 -> -6115
$ws.Cells.Item(79, 8).Value = 3848.229  # H79: 6529.407 -> 3848.229
$ws.Cells.Item(79, 9).Value = 3328.6191  # I79: 5199 -> 3328.6191
$ws.Cells.Item(79, 10).Value = 5485  # J79: 7312 -> 5485
$ws.Cells.Item(79, 11).Value = 3328.6191  # K79: 5199 -> 3328.6191
$ws.Cells.Item(79, 12).Value = 5485  # L79: 7312 -> 5485
$ws.Cells.Item(79, 13).Value = -2236.6191  # M79: -4107 -> -2236.6191
$ws.Cells.Item(79, 14).Value = -7669  # N79: -9496 -> -7669
$ws.Cells.Item(92, 8).Value = 3185.28  # H92: 3654.4092 -> 3185.28
$ws.Cells.Item(92, 9).Value = 1234.6666  # I92: 1464.85 -> 1234.6666
$ws.Cells.Item(92, 10).Value = 50000  # J92: 25550 -> 50000
$ws.Cells.Item(92, 11).Value = 1234.6666  # K92: 1464.85 -> 1234.6666
$ws.Cells.Item(92, 12).Value = 50000  # L92: 25550 -> 50000
$ws.Cells.Item(92, 13).Value = 13.33339999999998  # M92: -216.8499999999999 -> 13.33339999999998
$ws.Cells.Item(92, 14).Value = -52496  # N92: -28046 -> -52496
$ws.Cells.Item(98, 8).Value = 2290.4167  # H98: 1042.027 -> 2290.4167
$ws.Cells.Item(98, 9).Value = 2861  # I98: 884.6429000000001 -> 2861
$ws.Cells.Item(98, 10).Value = 1882.8572  # J98: 1531.6666 -> 1882.8572
$ws.Cells.Item(98, 11).Value = 2861  # K98: 884.6429000000001 -> 2861
$ws.Cells.Item(98, 12).Value = 1882.8572  # L98: 1531.6666 -> 1882.8572
$ws.Cells.Item(98, 13).Value = -1363  # M98: 613.3570999999999 -> -1363
$ws.Cells.Item(98, 14).Value = -4878.8572  # N98: -4527.6666 -> -4878.8572
$ws.Cells.Item(122, 8).Value = 2290.4167  # H122: 1042.027 -> 2290.4167
$ws.Cells.Item(122, 9).Value = 2861  # I122: 884.6429000000001 -> 2861
$ws.Cells.Item(122, 10).Value = 1882.8572  # J122: 1531.6666 -> 1882.8572
$ws.Cells.Item(122, 11).Value = 8583  # K122: 2653.9287 -> 8583
$ws.Cells.Item(122, 12).Value = 5648.571599999999  # L122: 4594.9998 -> 5648.571599999999
$ws.Cells.Item(122, 13).Value = -6133  # M122: -203.9287000000004 -> -6133
$ws.Cells.Item(122, 14).Value = -10548.5716  # N122: -9494.9998 -> -10548.5716
$ws.Cells.Item(127, 8).Value = 768977.7  # H127: 769056.3 -> 768977.7
$ws.Cells.Item(127, 9).Value = 450  # I127: 0 -> 450
$ws.Cells.Item(127, 10).Value = 897065.7  # J127: 769056.3 -> 897065.7
$ws.Cells.Item(127, 11).Value = 1350  # K127: 0 -> 1350
$ws.Cells.Item(127, 12).Value = 2691197.1  # L127: 2307168.9 -> 2691197.1
$ws.Cells.Item(127, 13).Value = 3610  # M127: None -> 3610
$ws.Cells.Item(127, 14).Value = -2701117.1  # N127: -2317088.9 -> -2701117.1
$ws.Cells.Item(137, 8).Value = 5408247  # H137: 6670088 -> 5408247
$ws.Cells.Item(137, 9).Value = 1869.0869  # I137: 2294.4707 -> 1869.0869
$ws.Cells.Item(137, 10).Value = 14290154  # J137: 15389510 -> 14290154
$ws.Cells.Item(137, 11).Value = 5607.2607  # K137: 6883.4121 -> 5607.2607
$ws.Cells.Item(137, 12).Value = 42870462  # L137: 46168530 -> 42870462
$ws.Cells.Item(137, 13).Value = -3057.2607  # M137: -4333.4121 -> -3057.2607
$ws.Cells.Item(137, 14).Value = -42875562  # N137: -46173630 -> -42875562
$ws.Cells.Item(138, 8).Value = 7578051  # H138: 5683748.5 -> 7578051
$ws.Cells.Item(138, 9).Value = 2036.9231  # I138: 1606.2354 -> 2036.9231
$ws.Cells.Item(138, 10).Value = 35717532  # J138: 25003032 -> 35717532
$ws.Cells.Item(138, 11).Value = 6110.7693  # K138: 4818.706200000001 -> 6110.7693
$ws.Cells.Item(138, 12).Value = 107152596  # L138: 75009096 -> 107152596
$ws.Cells.Item(138, 13).Value = -970.7692999999999  # M138: 321.2937999999995 -> -970.7692999999999
$ws.Cells.Item(138, 14).Value = -107162876  # N138: -75019376 -> -107162876
$ws.Cells.Item(141, 8).Value = 1706.0435  # H141: 1853.3636 -> 1706.0435
$ws.Cells.Item(141, 9).Value = 1157.7222  # I141: 1273.375 -> 1157.7222
$ws.Cells.Item(141, 10).Value = 3680  # J141: 3400 -> 3680
$ws.Cells.Item(141, 11).Value = 3473.1666  # K141: 3820.125 -> 3473.1666
$ws.Cells.Item(141, 12).Value = 11040  # L141: 10200 -> 11040
$ws.Cells.Item(141, 13).Value = 1706.8334  # M141: 1359.875 -> 1706.8334
$ws.Cells.Item(141, 14).Value = -21400  # N141: -20560 -> -21400

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(48, 8).Value = 245337  # H48: 289990 -> 245337
$ws.Cells.Item(48, 10).Value = 245337  # J48: 289990 -> 245337
$ws.Cells.Item(48, 12).Value = 245337  # L48: 289990 -> 245337
$ws.Cells.Item(48, 14).Value = -246105  # N48: -290758 -> -246105
$ws.Cells.Item(61, 8).Value = 9806070  # H61: 11366072 -> 9806070
$ws.Cells.Item(61, 9).Value = 13515376  # I61: 16668878 -> 13515376
$ws.Cells.Item(61, 10).Value = 2903  # J61: 2917.2856 -> 2903
$ws.Cells.Item(61, 11).Value = 13515376  # K61: 16668878 -> 13515376
$ws.Cells.Item(61, 12).Value = 2903  # L61: 2917.2856 -> 2903
$ws.Cells.Item(61, 13).Value = -13515164  # M61: -16668666 -> -13515164
$ws.Cells.Item(61, 14).Value = -3327  # N61: -3341.2856 -> -3327
$ws.Cells.Item(136, 8).Value = 9806070  # H136: 11366072 -> 9806070
$ws.Cells.Item(136, 9).Value = 13515376  # I136: 16668878 -> 13515376
$ws.Cells.Item(136, 10).Value = 2903  # J136: 2917.2856 -> 2903
$ws.Cells.Item(136, 11).Value = 40546128  # K136: 50006634 -> 40546128
$ws.Cells.Item(136, 12).Value = 8709  # L136: 8751.856800000001 -> 8709
$ws.Cells.Item(136, 13).Value = -40543578  # M136: -50004084 -> -40543578
$ws.Cells.Item(136, 14).Value = -13809  # N136: -13851.8568 -> -13809
$ws.Cells.Item(140, 8).Value = 65389  # H140: 46285.332 -> 65389
$ws.Cells.Item(140, 10).Value = 65389  # J140: 46285.332 -> 65389
$ws.Cells.Item(140, 12).Value = 65389  # L140: 46285.332 -> 65389
$ws.Cells.Item(140, 14).Value = -75749  # N140: -56645.332 -> -75749

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 0  # H16: 5004 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 5004 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 5004 -> 0
$ws.Cells.Item(16, 13).Value = ""  # M16: -4834 -> (removed)
$ws.Cells.Item(22, 8).Value = 453.1111  # H22: 415 -> 453.1111
$ws.Cells.Item(22, 9).Value = 486.85715  # I22: 415 -> 486.85715
$ws.Cells.Item(22, 10).Value = 335  # J22: 0 -> 335
$ws.Cells.Item(22, 11).Value = 486.85715  # K22: 415 -> 486.85715
$ws.Cells.Item(22, 12).Value = 335  # L22: 0 -> 335
$ws.Cells.Item(22, 13).Value = -313.85715  # M22: -242 -> -313.85715
$ws.Cells.Item(22, 14).Value = -681  # N22: None -> -681

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7940630.5  # H31: 6413678 -> 7940630.5
$ws.Cells.Item(31, 9).Value = 8009.0557  # I31: 5015.5 -> 8009.0557
$ws.Cells.Item(31, 10).Value = 13890096  # J31: 15152764 -> 13890096
$ws.Cells.Item(31, 11).Value = 8009.0557  # K31: 5015.5 -> 8009.0557
$ws.Cells.Item(31, 12).Value = 13890096  # L31: 15152764 -> 13890096
$ws.Cells.Item(31, 13).Value = -7714.0557  # M31: -4720.5 -> -7714.0557
$ws.Cells.Item(31, 14).Value = -13890686  # N31: -15153354 -> -13890686
$ws.Cells.Item(34, 8).Value = 7940630.5  # H34: 6413678 -> 7940630.5
$ws.Cells.Item(34, 9).Value = 8009.0557  # I34: 5015.5 -> 8009.0557
$ws.Cells.Item(34, 10).Value = 13890096  # J34: 15152764 -> 13890096
$ws.Cells.Item(34, 11).Value = 8009.0557  # K34: 5015.5 -> 8009.0557
$ws.Cells.Item(34, 12).Value = 13890096  # L34: 15152764 -> 13890096
$ws.Cells.Item(34, 13).Value = -7807.0557  # M34: -4813.5 -> -7807.0557
$ws.Cells.Item(34, 14).Value = -13890500  # N34: -15153168 -> -13890500
$ws.Cells.Item(86, 8).Value = 1732.0952  # H86: 2370.8 -> 1732.0952
$ws.Cells.Item(86, 9).Value = 1858.8  # I86: 2750 -> 1858.8
$ws.Cells.Item(86, 10).Value = 1616.909  # J86: 2118 -> 1616.909
$ws.Cells.Item(86, 11).Value = 1858.8  # K86: 2750 -> 1858.8
$ws.Cells.Item(86, 12).Value = 1616.909  # L86: 2118 -> 1616.909
$ws.Cells.Item(86, 13).Value = -735.8  # M86: -1627 -> -735.8
$ws.Cells.Item(86, 14).Value = -3862.909  # N86: -4364 -> -3862.909
$ws.Cells.Item(89, 8).Value = 1732.0952  # H89: 2370.8 -> 1732.0952
$ws.Cells.Item(89, 9).Value = 1858.8  # I89: 2750 -> 1858.8
$ws.Cells.Item(89, 10).Value = 1616.909  # J89: 2118 -> 1616.909
$ws.Cells.Item(89, 11).Value = 9294  # K89: 13750 -> 9294
$ws.Cells.Item(89, 12).Value = 8084.545  # L89: 10590 -> 8084.545
$ws.Cells.Item(89, 13).Value = -3678  # M89: -8134 -> -3678
$ws.Cells.Item(89, 14).Value = -19316.545  # N89: -21822 -> -19316.545

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 717.0345  # H68: 689.2135 -> 717.0345
$ws.Cells.Item(68, 9).Value = 540.2785  # I68: 529.7560999999999 -> 540.2785
$ws.Cells.Item(68, 10).Value = 2462.5  # J68: 2557.1428 -> 2462.5
$ws.Cells.Item(68, 11).Value = 1620.8355  # K68: 1589.2683 -> 1620.8355
$ws.Cells.Item(68, 12).Value = 7387.5  # L68: 7671.428400000001 -> 7387.5
$ws.Cells.Item(68, 13).Value = -809.8355000000001  # M68: -778.2682999999997 -> -809.8355000000001
$ws.Cells.Item(68, 14).Value = -9009.5  # N68: -9293.428400000001 -> -9009.5
$ws.Cells.Item(71, 8).Value = 717.0345  # H71: 689.2135 -> 717.0345
$ws.Cells.Item(71, 9).Value = 540.2785  # I71: 529.7560999999999 -> 540.2785
$ws.Cells.Item(71, 10).Value = 2462.5  # J71: 2557.1428 -> 2462.5
$ws.Cells.Item(71, 11).Value = 4862.5065  # K71: 4767.804899999999 -> 4862.5065
$ws.Cells.Item(71, 12).Value = 22162.5  # L71: 23014.2852 -> 22162.5
$ws.Cells.Item(71, 13).Value = -806.5065000000004  # M71: -711.8048999999992 -> -806.5065000000004
$ws.Cells.Item(71, 14).Value = -30274.5  # N71: -31126.2852 -> -30274.5
$ws.Cells.Item(102, 8).Value = 3815.75  # H102: 8112.375 -> 3815.75
$ws.Cells.Item(102, 9).Value = 2542  # I102: 0 -> 2542
$ws.Cells.Item(102, 10).Value = 4580  # J102: 8112.375 -> 4580
$ws.Cells.Item(102, 11).Value = 7626  # K102: 0 -> 7626
$ws.Cells.Item(102, 12).Value = 13740  # L102: 24337.125 -> 13740
$ws.Cells.Item(102, 13).Value = -5192  # M102: None -> -5192
$ws.Cells.Item(102, 14).Value = -18608  # N102: -29205.125 -> -18608
$ws.Cells.Item(107, 8).Value = 968.0571  # H107: 1880.8572 -> 968.0571
$ws.Cells.Item(107, 9).Value = 407  # I107: 0 -> 407
$ws.Cells.Item(107, 10).Value = 1917.5385  # J107: 1880.8572 -> 1917.5385
$ws.Cells.Item(107, 11).Value = 1221  # K107: 0 -> 1221
$ws.Cells.Item(107, 12).Value = 5752.6155  # L107: 5642.571599999999 -> 5752.6155
$ws.Cells.Item(107, 13).Value = 699  # M107: None -> 699
$ws.Cells.Item(107, 14).Value = -9592.6155  # N107: -9482.571599999999 -> -9592.6155
$ws.Cells.Item(114, 8).Value = 1402.875  # H114: 1283.4286 -> 1402.875
$ws.Cells.Item(114, 9).Value = 453.25  # I114: 314 -> 453.25
$ws.Cells.Item(114, 10).Value = 2352.5  # J114: 3028.4 -> 2352.5
$ws.Cells.Item(114, 11).Value = 1359.75  # K114: 942 -> 1359.75
$ws.Cells.Item(114, 12).Value = 7057.5  # L114: 9085.200000000001 -> 7057.5
$ws.Cells.Item(114, 13).Value = 1894.25  # M114: 2312 -> 1894.25
$ws.Cells.Item(114, 14).Value = -13565.5  # N114: -15593.2 -> -13565.5
$ws.Cells.Item(117, 8).Value = 1021.125  # H117: 1027.375 -> 1021.125
$ws.Cells.Item(117, 10).Value = 1154.6923  # J117: 1162.3846 -> 1154.6923
$ws.Cells.Item(117, 12).Value = 3464.0769  # L117: 3487.1538 -> 3464.0769
$ws.Cells.Item(117, 14).Value = -10348.0769  # N117: -10371.1538 -> -10348.0769
$ws.Cells.Item(121, 8).Value = 1303.0667  # H121: 1207.1538 -> 1303.0667
$ws.Cells.Item(121, 9).Value = 112.25  # I121: 261.14285 -> 112.25
$ws.Cells.Item(121, 10).Value = 1736.091  # J121: 2310.8333 -> 1736.091
$ws.Cells.Item(121, 11).Value = 336.75  # K121: 783.4285500000001 -> 336.75
$ws.Cells.Item(121, 12).Value = 5208.272999999999  # L121: 6932.499899999999 -> 5208.272999999999
$ws.Cells.Item(121, 13).Value = 973.25  # M121: 526.5714499999999 -> 973.25
$ws.Cells.Item(121, 14).Value = -7828.272999999999  # N121: -9552.499899999999 -> -7828.272999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 23362.355  # H70: 24772.932 -> 23362.355
$ws.Cells.Item(70, 9).Value = 77772.125  # I70: 69719.664 -> 77772.125
$ws.Cells.Item(70, 10).Value = 4437.2173  # J70: 4546.9 -> 4437.2173
$ws.Cells.Item(70, 11).Value = 77772.125  # K70: 69719.664 -> 77772.125
$ws.Cells.Item(70, 12).Value = 4437.2173  # L70: 4546.9 -> 4437.2173
$ws.Cells.Item(70, 13).Value = -77502.125  # M70: -69449.664 -> -77502.125
$ws.Cells.Item(70, 14).Value = -4977.2173  # N70: -5086.9 -> -4977.2173
$ws.Cells.Item(73, 8).Value = 23362.355  # H73: 24772.932 -> 23362.355
$ws.Cells.Item(73, 9).Value = 77772.125  # I73: 69719.664 -> 77772.125
$ws.Cells.Item(73, 10).Value = 4437.2173  # J73: 4546.9 -> 4437.2173
$ws.Cells.Item(73, 11).Value = 77772.125  # K73: 69719.664 -> 77772.125
$ws.Cells.Item(73, 12).Value = 4437.2173  # L73: 4546.9 -> 4437.2173
$ws.Cells.Item(73, 13).Value = -76836.125  # M73: -68783.664 -> -76836.125
$ws.Cells.Item(73, 14).Value = -6309.2173  # N73: -6418.9 -> -6309.2173
$ws.Cells.Item(126, 8).Value = 5700  # H126: 4306.846 -> 5700
$ws.Cells.Item(126, 9).Value = 3000  # I126: 1997.8 -> 3000
$ws.Cells.Item(126, 10).Value = 6240  # J126: 5750 -> 6240
$ws.Cells.Item(126, 11).Value = 9000  # K126: 5993.4 -> 9000
$ws.Cells.Item(126, 12).Value = 18720  # L126: 17250 -> 18720
$ws.Cells.Item(126, 13).Value = -6530  # M126: -3523.4 -> -6530
$ws.Cells.Item(126, 14).Value = -23660  # N126: -22190 -> -23660

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 0  # H57: 50000 -> 0
$ws.Cells.Item(57, 10).Value = 0  # J57: 50000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 50000 -> 0
$ws.Cells.Item(57, 14).Value = ""  # N57: -51508 -> (removed)
$ws.Cells.Item(126, 8).Value = 4163.2  # H126: 3084.9333 -> 4163.2
$ws.Cells.Item(126, 9).Value = 2034.4615  # I126: 1932.2354 -> 2034.4615
$ws.Cells.Item(126, 10).Value = 18000  # J126: 4592.3076 -> 18000
$ws.Cells.Item(126, 11).Value = 6103.3845  # K126: 5796.706200000001 -> 6103.3845
$ws.Cells.Item(126, 12).Value = 54000  # L126: 13776.9228 -> 54000
$ws.Cells.Item(126, 13).Value = -3633.3845  # M126: -3326.706200000001 -> -3633.3845
$ws.Cells.Item(126, 14).Value = -58940  # N126: -18716.9228 -> -58940
$ws.Cells.Item(132, 8).Value = 3348.3958  # H132: 4600.9707 -> 3348.3958
$ws.Cells.Item(132, 9).Value = 3666.4211  # I132: 5279.423 -> 3666.4211
$ws.Cells.Item(132, 10).Value = 2139.9  # J132: 2396 -> 2139.9
$ws.Cells.Item(132, 11).Value = 10999.2633  # K132: 15838.269 -> 10999.2633
$ws.Cells.Item(132, 12).Value = 6419.700000000001  # L132: 7188 -> 6419.700000000001
$ws.Cells.Item(132, 13).Value = -8469.263300000001  # M132: -13308.269 -> -8469.263300000001
$ws.Cells.Item(132, 14).Value = -11479.7  # N132: -12248 -> -11479.7
